$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Split the old generic "STR24" STR_ID into "STR24A" / "STR24B"
#    depending on the row (rows 34-43)
# -----------------------------------------------------------------
$str24a = @(34, 36, 42, 43)
$str24b = @(35, 37, 38, 39, 40, 41)

foreach ($r in $str24a) {
    $ws.Cells.Item($r, 1).Value = "STR24A"
}
foreach ($r in $str24b) {
    $ws.Cells.Item($r, 1).Value = "STR24B"
}

# -----------------------------------------------------------------
# 2. Add the new "Analyse" column (U) = 0 for every existing data row
#    (rows 2-43)
# -----------------------------------------------------------------
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 21).Value = 0
}

# -----------------------------------------------------------------
# 3. Add missing "Passed_QC" column (T) = 0 for rows 31-43 (these rows
#    did not have a Passed_QC value yet)
# -----------------------------------------------------------------
for ($r = 31; $r -le 43; $r++) {
    $ws.Cells.Item($r, 20).Value = 0
}

# -----------------------------------------------------------------
# 5. Append four new data rows (44-47) for the new screen "STR27"
# -----------------------------------------------------------------

# Carry over the date number-format (style) used by the other date
# columns (B/C) onto the new rows first.
$ws.Range("B2:C2").Copy()
$ws.Range("B44:C47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 44; A = "STR27"; B = 44939; C = 44944; D = "FullscreenV2.xlsx"; E = "OPT0034"; G = "OPT0034"; H = 0; I = 1; J = 3; L = 3;  M = 0; N = 3; O = "3D"; Q = 0; T = 0; U = 1 },
    @{ Row = 45; A = "STR27"; B = 44939; C = 44944; D = "FullscreenV2.xlsx"; E = "OPT0413"; G = "OPT0413"; H = 0; I = 1; J = 3; L = 1;  M = 0; N = 3; O = "3D"; Q = 0; T = 0; U = 1 },
    @{ Row = 46; A = "STR27"; B = 44939; C = 44944; D = "FullscreenV2.xlsx"; E = "RAS11";   G = "RAS11";   H = 1; I = 0; J = 3; L = 7;  M = 0; N = 3; O = "3D"; Q = 0; T = 0; U = 1 },
    @{ Row = 47; A = "STR27"; B = 44939; C = 44944; D = "FullscreenV2.xlsx"; E = "RAS27";   G = "RAS27";   H = 1; I = 0; J = 3; L = 5;  M = 0; N = 3; O = "3D"; Q = 0; T = 0; U = 1 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = $row.B
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
}

# -----------------------------------------------------------------
# 5b. Add new "Analyse" header in U1 (added last so it lands after
#     the other newly introduced shared strings)
# -----------------------------------------------------------------
$ws.Cells.Item(1, 21).Value = "Analyse"

# -----------------------------------------------------------------
# 6. Restore the frozen-pane scroll position / selection, matching
#    where the user ended up after the edits.
# -----------------------------------------------------------------
$ws.Range("D25").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("Q48:U48").Select()
